# Regenerate orders with updated distance/sizes:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# These substitutions are applied to every cell value in the used range
# (condition labels, filenames, Distance/Size lookup columns, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

$used.Replace("D64", "D69", $false, $false, $false) | Out-Null
$used.Replace("D51", "D55", $false, $false, $false) | Out-Null
$used.Replace("D80", "D86", $false, $false, $false) | Out-Null
$used.Replace("S30", "S31", $false, $false, $false) | Out-Null
